$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card8")

# 1. Add new header cell M1 "Event " with the same style as the other header cells (A1:L1)
$ws.Range("A1").Copy()
$ws.Range("M1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M1").Value = "Event "

# 2. Populate M2:M13 as present-but-empty text cells (matches the blank "Event" column body)
for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 13)
    $cell.Value = "'"
    $cell.Style = "Normal"
}

# 3. Row 8 (F8:K8) previously blank cells now hold literal "nan" text
$ws.Range("F8:K8").Value = "nan"
